$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Table S3")
$ws.Range("D3").Value = 'ns (p = .05)'
$ws.Range("E3").Value = 'r = .19'
$ws.Range("D4").Value = 'ns (p = .10)'
$ws.Range("E4").Value = 'V = .19'
$ws.Range("D5").Value = 'ns (p = .05)'
$ws.Range("E5").Value = 'V = .23'
$ws.Range("D6").Value = 'p = .03'
$ws.Range("E6").Value = 'r = .2'
$ws.Range("D7").Value = 'ns (p = .09)'
$ws.Range("E7").Value = 'r = .16'
$ws.Range("D8").Value = 'ns (p = .05)'
$ws.Range("E8").Value = 'r = .19'
$ws = $wb.Worksheets.Item("Table S4")
$ws.Range("F3").Value = 'ns (p = .55)'
$ws.Range("G3").Value = 'V = .13'
$ws.Range("F4").Value = 'ns (p = .70)'
$ws.Range("G4").Value = 'η² = .018'
$ws.Range("F5").Value = 'ns (p = .98)'
$ws.Range("G5").Value = 'V = .025'
$ws.Range("F6").Value = 'ns (p = .65)'
$ws.Range("G6").Value = 'V = .11'
$ws.Range("F7").Value = 'ns (p = .68)'
$ws.Range("G7").Value = 'V = .1'
$ws.Range("F8").Value = 'ns (p = .28)'
$ws.Range("G8").Value = 'V = .19'
$ws.Range("F9").Value = 'ns (p = .18)'
$ws.Range("G9").Value = 'V = .21'
$ws.Range("F10").Value = 'ns (p = .13)'
$ws.Range("G10").Value = 'V = .23'
$ws.Range("F11").Value = 'ns (p = .55)'
$ws.Range("G11").Value = 'V = .13'
$ws.Range("F12").Value = 'ns (p = .49)'
$ws.Range("G12").Value = 'V = .14'
$ws.Range("F13").Value = 'ns (p = .68)'
$ws.Range("G13").Value = 'V = .1'
$ws.Range("F14").Value = 'ns (p = .54)'
$ws.Range("G14").Value = 'η² = .011'
$ws.Range("F15").Value = 'ns (p = .29)'
$ws.Range("G15").Value = 'V = .18'
$ws.Range("F16").Value = 'ns (p = .21)'
$ws.Range("G16").Value = 'V = .2'
$ws.Range("F17").Value = 'p < .001'
$ws.Range("G17").Value = 'V = .52'
$ws.Range("F18").Value = 'p < .001'
$ws.Range("G18").Value = 'η² = .35'
$ws.Range("F19").Value = 'p = .002'
$ws.Range("G19").Value = 'V = .41'
$ws.Range("F20").Value = 'ns (p = .11)'
$ws.Range("G20").Value = 'η² = .034'
$ws.Range("F21").Value = 'ns (p = .85)'
$ws.Range("G21").Value = 'η² = .024'
$ws.Range("F22").Value = 'p < .001'
$ws.Range("G22").Value = 'V = .61'
$ws = $wb.Worksheets.Item("Table S5")
$ws.Range("F3").Value = 'ns (p = .28)'
$ws.Range("G3").Value = 'V = .19'
$ws.Range("F4").Value = 'p = .003'
$ws.Range("G4").Value = 'η² = .14'
$ws.Range("F5").Value = 'ns (p = .19)'
$ws.Range("G5").Value = 'η² = .018'
$ws.Range("F6").Value = 'ns (p = .82)'
$ws.Range("G6").Value = 'η² = .022'
$ws.Range("F7").Value = 'ns (p = .49)'
$ws.Range("G7").Value = 'η² = .008'
$ws.Range("F8").Value = 'ns (p = .28)'
$ws.Range("G8").Value = 'η² = .0078'
$ws.Range("F9").Value = 'ns (p = .18)'
$ws.Range("G9").Value = 'η² = .021'
$ws.Range("F10").Value = 'ns (p = .13)'
$ws.Range("G10").Value = 'η² = .029'
$ws.Range("F11").Value = 'ns (p = .80)'
$ws.Range("G11").Value = 'η² = .022'
$ws.Range("F12").Value = 'ns (p = .11)'
$ws.Range("G12").Value = 'η² = .035'
$ws = $wb.Worksheets.Item("Table S6")
$ws.Range("D3").Value = 'ns (p = .09)'
$ws.Range("E3").Value = 'V = .22'
$ws.Range("D4").Value = 'p = .03'
$ws.Range("E4").Value = 'V = .3'
$ws.Range("D5").Value = 'ns (p = .07)'
$ws.Range("E5").Value = 'V = .24'
$ws.Range("D6").Value = 'p = .002'
$ws.Range("E6").Value = 'r = .37'
$ws.Range("D7").Value = 'p = .01'
$ws.Range("E7").Value = 'V = .34'
$ws.Range("D8").Value = 'p = .04'
$ws.Range("E8").Value = 'V = .27'
$ws.Range("D9").Value = 'p = .02'
$ws.Range("E9").Value = 'r = .27'
$ws.Range("D10").Value = 'ns (p = .09)'
$ws.Range("E10").Value = 'V = .23'
$ws.Range("D11").Value = 'p < .001'
$ws.Range("E11").Value = 'V = .51'
$ws.Range("D12").Value = 'p < .001'
$ws.Range("E12").Value = 'r = .49'
$ws.Range("D13").Value = 'p < .001'
$ws.Range("E13").Value = 'r = .51'
$ws.Range("D14").Value = 'p = .02'
$ws.Range("E14").Value = 'r = .26'
$ws.Range("D15").Value = 'ns (p = .26)'
$ws.Range("E15").Value = 'r = .13'
$ws.Range("D16").Value = 'ns (p > .99)'
$ws.Range("E16").Value = 'r = .0013'
$ws.Range("D17").Value = 'ns (p = .26)'
$ws.Range("E17").Value = 'r = .13'
$ws.Range("D18").Value = 'ns (p = .68)'
$ws.Range("E18").Value = 'r = .049'
$ws.Range("D19").Value = 'ns (p = .47)'
$ws.Range("E19").Value = 'r = .085'
$ws.Range("D20").Value = 'ns (p = .26)'
$ws.Range("E20").Value = 'r = .13'
$ws.Range("D21").Value = 'ns (p = .77)'
$ws.Range("E21").Value = 'r = .035'
$ws.Range("D22").Value = 'ns (p = .72)'
$ws.Range("E22").Value = 'r = .042'
$ws.Range("D23").Value = 'ns (p = .66)'
$ws.Range("E23").Value = 'r = .053'
$ws.Range("D24").Value = 'ns (p = .19)'
$ws.Range("E24").Value = 'r = .15'
$ws.Range("D25").Value = 'ns (p = .78)'
$ws.Range("E25").Value = 'r = .034'
$ws = $wb.Worksheets.Item("Table S7")
$ws.Range("D3").Value = 'p = .01'
$ws.Range("E3").Value = 'V = .33'
$ws.Range("D4").Value = 'p = .002'
$ws.Range("E4").Value = 'r = .37'
$ws.Range("D5").Value = 'p = .001'
$ws.Range("E5").Value = 'V = .42'
$ws.Range("D6").Value = 'p = .02'
$ws.Range("E6").Value = 'r = .27'
$ws.Range("D7").Value = 'ns (p = .08)'
$ws.Range("E7").Value = 'V = .24'
$ws.Range("D8").Value = 'p = .03'
$ws.Range("E8").Value = 'V = .28'
$ws.Range("D9").Value = 'p = .005'
$ws.Range("E9").Value = 'r = .33'
$ws.Range("D10").Value = 'p = .03'
$ws.Range("E10").Value = 'r = .26'
$ws.Range("D11").Value = 'p = .03'
$ws.Range("E11").Value = 'r = .25'
$ws.Range("D12").Value = 'p = .006'
$ws.Range("E12").Value = 'r = .32'
$ws.Range("D13").Value = 'ns (p = .05)'
$ws.Range("E13").Value = 'r = .23'
$ws = $wb.Worksheets.Item("Table S8")
$ws.Range("E3").Value = 'ns (p = .09)'
$ws.Range("F3").Value = 'V = .26'
$ws.Range("E4").Value = 'ns (p = .09)'
$ws.Range("F4").Value = 'V = .26'
$ws.Range("E5").Value = 'ns (p = .09)'
$ws.Range("F5").Value = 'V = .25'
$ws.Range("E6").Value = 'p = .048'
$ws.Range("F6").Value = 'η² = .058'
$ws.Range("E7").Value = 'p = .02'
$ws.Range("F7").Value = 'V = .33'
$ws.Range("E8").Value = 'p < .001'
$ws.Range("F8").Value = 'η² = .31'
$ws.Range("E9").Value = 'p < .001'
$ws.Range("F9").Value = 'V = .51'
$ws.Range("E10").Value = 'p = .006'
$ws.Range("F10").Value = 'V = .37'
$ws.Range("E11").Value = 'p = .008'
$ws.Range("F11").Value = 'V = .36'
$ws.Range("E12").Value = 'p = .02'
$ws.Range("F12").Value = 'V = .32'
$ws.Range("E13").Value = 'p < .001'
$ws.Range("F13").Value = 'V = .57'
$ws.Range("E14").Value = 'p < .001'
$ws.Range("F14").Value = 'η² = .34'
$ws.Range("E15").Value = 'p < .001'
$ws.Range("F15").Value = 'V = .57'
$ws.Range("E16").Value = 'ns (p = .10)'
$ws.Range("F16").Value = 'V = .25'
$ws.Range("E17").Value = 'p = .02'
$ws.Range("F17").Value = 'η² = .09'
$ws.Range("E18").Value = 'p = .01'
$ws.Range("F18").Value = 'V = .35'
